$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header values first (AD1, AE1, AF1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from the existing
# header cell A1 onto the three new header cells without disturbing the
# values we just set.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record (Wins/Losses/Ties) for every data row (2-48).
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 70
    $ws.Cells.Item($r, 31).Value = 92
    $ws.Cells.Item($r, 32).Value = 0
}
